$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ------------------------------------------------------------------
# 1. "Ativação:" date changes from 01/01/2012 to 01/01/2023.
#    This value is shared by B8/C8 (Ativação row) and B15/C15 (which,
#    per the original workbook, also display the same date string).
#    We build the replacement text in a scratch cell first (forcing a
#    Text number format so Excel doesn't reinterpret "01/01/2023" as a
#    real date serial), then paste only the VALUE into each target
#    cell so their existing formatting (style) is left untouched.
# ------------------------------------------------------------------
$ws.Range("Z1").NumberFormat = "@"
$ws.Range("Z1").Value = "01/01/2023"
$ws.Range("Z1").Copy()

$ws.Range("B8").PasteSpecial(-4163)   # xlPasteValues
$ws.Range("C8").PasteSpecial(-4163)
$ws.Range("B15").PasteSpecial(-4163)
$ws.Range("C15").PasteSpecial(-4163)

$ws.Range("Z1").Clear()
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 2. Row 11 "Objectives:" - add the English objectives paragraph to
#    B11/C11 (columns previously empty on this row).
# ------------------------------------------------------------------
$ws.Range("B11").Value = "Provide the student with knowledge of the main techniques of physical and chemical characterization of materials."
$ws.Range("C11").Value = "Provide the student with knowledge of the main techniques of physical and chemical characterization of materials."
$ws.Range("B13").Copy()
$ws.Range("B11").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C13").Copy()
$ws.Range("C11").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 3. Row 14 "Short syllabus:" - add the English short syllabus text to
#    B14/C14.
# ------------------------------------------------------------------
$ws.Range("B14").Value = "Granulometric and surface analysis. Microstructural analyses. Thermal analysis. Rheometry."
$ws.Range("C14").Value = "Granulometric and surface analysis. Microstructural analyses. Thermal analysis. Rheometry."
$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ------------------------------------------------------------------
# 4. Row 16 "Syllabus:" - add the English syllabus text to B16/C16.
# ------------------------------------------------------------------
$ws.Range("B16").Value = "Grain size analysis. BET adsorption, porosity and pycnometry.Microstructural analysis: X-ray diffraction, Laue figure; X-ray scattering (SAXS). Electron diffraction. Optical Microscopy. Electron microscopy, X-ray microanalysis (EDX and WDX).Thermal analysis: Differential thermal analysis (DTA), differential scanning calorimetry (DSC) and thermogravimetry (TGA).Rheometry of liquids, solutions and pastes."
$ws.Range("C16").Value = "Grain size analysis. BET adsorption, porosity and pycnometry.Microstructural analysis: X-ray diffraction, Laue figure; X-ray scattering (SAXS). Electron diffraction. Optical Microscopy. Electron microscopy, X-ray microanalysis (EDX and WDX).Thermal analysis: Differential thermal analysis (DTA), differential scanning calorimetry (DSC) and thermogravimetry (TGA).Rheometry of liquids, solutions and pastes."
$ws.Range("B13").Copy()
$ws.Range("B16").PasteSpecial(-4122)
$ws.Range("C13").Copy()
$ws.Range("C16").PasteSpecial(-4122)
$excel.CutCopyMode = $false

$wb.Save()
